$d = $word.ActiveDocument

# Append a new run ", 心情也很好" (full-width Chinese comma) right after the
# existing "今天天气很好" text at the very end of the document, matching the
# size formatting (22pt == half-points 44) used throughout the document.

$range = $d.Content
$range.Collapse(0)  # wdCollapseEnd

$range.InsertAfter("，心情也很好")
$range.Font.Size = 22

$d.Save()
